$p = $ppt.ActivePresentation

# --- Slide 1 (existing title slide): set title + subtitle text ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = 'Organizing Oppias Learning Content'
$s1.Shapes.Item(2).TextFrame.TextRange.Text = 'Authors: Ben Henning, Sean Lip, Tony Jiang'

# --- Slide 2: Related References ---
$s = $p.Slides.Add(2, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Related References'
$body = @'
Functional Requirements: link
Overall UX and page URLs: 
Skills breakdown: 
Questions breakdown: 
Topics/stories breakdown: 
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 3: 1 Overview ---
$s = $p.Slides.Add(3, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = '1 Overview'

# --- Slide 4: 1.1 Background ---
$s = $p.Slides.Add(4, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = '1.1 Background'
$body = @'

Oppia currently primarily uses both collections and explorations to co-organize all learning content on Oppia. These activities, separately, are organized into subject categories to be displayed on the Oppia library page. There is currently work being done to introduce questions and refresher explorations as additional learning-related activities (see here for a much more in-depth background).

This proposal aims to:
Establish a clear, high-level goal used when organizing each aspect of Oppias learning content (section 1.2)
Enumerate expected user scenarios that demonstrate issues that require a solution to our goal (section 1.3)
Propose a structure that helps lead us toward reaching this goal and solving the specific user problems (section 2)
Briefly consider technical implications of this structure (section 3)
Provide a pathway to future work that may continue helping us approach the goal (section 3s open questions)
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 5: 1.2 Goal ---
$s = $p.Slides.Add(5, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = '1.2 Goal'
$body = @'
The goal of this taxonomy is to allow learners to jump in anywhere within Oppias educational landscape and be gently guided to where they should really start. This course correction is a central behavior which allows us to focus on one core piece of functionality while facilitating solutions to several expected user scenarios (discussed in section 1.3).
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 6: 1.3 Expected User Scenarios ---
$s = $p.Slides.Add(6, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = '1.3 Expected User Scenarios'
$body = @'
Following are expected personas of students using Oppia, though they are intentionally not fully fleshed out to keep their application in this document simple:
Primary: student using Oppia inside classroom environment
Primary: student using Oppia at home as a teacher-suggested supplement
Primary: student using Oppia at home as a personal supplement to an enrolled class
Secondary: student using Oppia to study content unrelated to their enrolled classes
Secondary: non-student using Oppia to study anything
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 7: 2 Proposed Structure ---
$s = $p.Slides.Add(7, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = '2 Proposed Structure'
$body = @'
The following paragraphs define a high-level taxonomy of constructs in Oppia.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 8: Subject ---
$s = $p.Slides.Add(8, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Subject'
$body = @'
Subjects describe the entire universe of knowledge related to a specific subject area (e.g. Mathematics or Physics). Subjects contain one or more topics. There are no firm guidelines yet on what a subject is, or how subjects should be split into topics. Instead, we will start by having a tight set of admins create a few subjects and iterate until the Oppia Foundation has formalized clear boundaries between subjects. We will then use these as examples for future subjects.

Learners viewing a subject will see a list of topics that they can learn, with search functionality. The library page will show individual subject pages once they are available, instead of categories of explorations and collections. In the meantime, topics will be listed in admin-defined groups.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 9: Topic ---
$s = $p.Slides.Add(9, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Topic'
$body = @'
A topic is a high-level concepts related to a certain subject (e.g. differentiation or Newtons laws of motion). There are no firm guidelines yet on what a topic is or how certain educational material should be split across multiple topics. Instead, we will start by having a tight set of admins create a few topics for a given subject, and iterate until the Oppia Foundation has formalized clear boundaries between topics. We will then use these as examples for future topics.

When learning a topic, a learner has access to a complete package of learning, including  narratives (stories) that teach the topic in detail. Each topic may have one or more story arcs associated with it. The learner has the option, at any time, to practice what theyve learned in that topic using practice sessions.

In the frontend, topics have a landing page which provides the learner with some initial context on what they will be learning, or an optional reference to the stories being taught. The landing page provides more insight into what the topic means, and aims to pique students interest in learning the topic.

Users may also fast-track the skills taught by a topic by opening the concept cards for each skill and practicing them directly.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 10: Story (Narrative) ---
$s = $p.Slides.Add(10, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Story (Narrative)'
$body = @'
A story guides the learner through a journey of lessons to teach a specific topic in full. Each story helps provide additional context when learning the skills of the topic.

Lessons will be prefixed with a ~3 question, non-context-specific pre-test which will evaluate whether the learner has the expected skills to begin that lesson. This test is generated by the platform and allows learners to begin studying a topic without necessarily having all of the prerequisites needed to begin studying. If the learner cannot complete any of these questions, they will be shown a concept card reviewing that material in more detail.

Every 2-3 lessons the learner will complete a generated, non-context-specific review test of ~10 questions to double check that they understand the content covered in prior lessons in the story before moving onward. See the section on review tests below for more detail.

The learner is discouraged (but not disallowed) from skipping ahead in the story without first completing an earlier part. This is done by hiding later lessons in the story using a growing world map (where the metaphor is continued using actual animations to help improve story immersion; think Sid Meiers Civilization for world hiding and movement).

Each topic will include a list of canonical stories. These stories, taken together, should teach all the skill of the given topic (but an individual story does not need to teach every skill). Its up to the creator to decide how large or broad a given story should be.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 11: Lesson ---
$s = $p.Slides.Add(11, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Lesson'
$body = @'
A lesson is an exploration which teaches one or more skills in depth by putting them within the context of a story. In general, each lesson is expected to teach a specific, coherent, and self-contained idea that is coherently described in the context of a topic (e.g. adding fractions, or the chain rule of differentiation).

Lessons will include questions that take place within the context of the story. These questions focus on pinning down misconceptions early. They are not taken from the question bank, but they may still be associated with a given skill. If Oppia detects that a student is struggling with a certain prerequisite skill when answering a question, the lesson may reinforce that prerequisite skill by surfacing the relevant concept card and providing a 3-question skill test that must be completed before the learner can continue.

Lessons are expected to include questions at the end to review the new content taught within that exploration before moving on to the next lesson, where wrong answers lead to the learner being guided to an earlier part of the lesson to refresh their understanding of that topic. 

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 12: Concept Card ---
$s = $p.Slides.Add(12, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Concept Card'
$body = @'
A concept card is a detailed overview of the material for a given skill (e.g. an explanation and learner-selected, variable number of worked examples (pre-solved questions & answers). Learners can select more worked examples as desired. The learner may also, from the concept card, begin a practice session.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 13: Pre-Test ---
$s = $p.Slides.Add(13, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Pre-Test'
$body = @'
A pre-test is a variably long set of questions that has been randomly created from the list of prerequisite skills specified by a lesson. Each question must be answered correctly (possibly with the assistance of hints or a solution) before continuing. Since pre-tests are intended to gate access to a lesson until sufficient expertise in the prerequisite skills is demonstrated, a specific score needs to be earned before the pre-test is considered passed. See the Score section for details on how scores work.



Numerical scores are not shown to the learner in pre-tests, but we may show a gamified progress bar to indicate how much work is left to complete in the pre-test.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 14: Review Test ---
$s = $p.Slides.Add(14, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Review Test'
$body = @'
Review tests are meant to help the learner refresh their understanding of skills theyve already studied within the topic before continuing to new concepts. A review test is a variable-length set of questions that has been randomly generated from the list of acquired skills from all previously completed lessons in the story arc for that topic (skills acquired from other story arcs in the lesson are not considered).

Similar to pre-tests, review tests require a certain score to be achieved before passing them. See the Pre-Test section for details on how scoring works for tests. Oppia will detect struggling for specific questions and show a concept card to the learner, if needed (see the Question section for more details).
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 15: Skill Test ---
$s = $p.Slides.Add(15, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Skill Test'
$body = @'
Skill tests are similar to pre-tests in that they must be passed in order to continue with a lesson, though skill tests may show up mid-lesson. They lack context, are generated by Oppia, and always focus on specifically one skill. Lessons require prerequisite skills in order to complete them, and if Oppia detects a learner is struggling on one of these skills it will pause the lesson in order to show the learner that skills concept card and require the learner to pass a skill test before proceeding.

Skill tests will appear as new lesson cards (natively within the learner view) rather than redirect the learner. These cards may be specially marked to note they are part of a temporary skill review test.

The skill test is around 3 questions and is meant to ensure the learner has reviewed the concept card and worked examples in detail to sufficiently demonstrate they understand the topic before proceeding with their current lesson.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 16: Practice Session ---
$s = $p.Slides.Add(16, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Practice Session'
$body = @'
Practice sessions are variable-length sets of questions randomly selected from a list of skills. If the practice session is started from a concept card, the only skill used to generate the questions is the skill associated with that concept card (the learner has no choice in this context). If the practice session is started from a topic, the learner may select which skills from that topic they would like to practice. Long-term, we may introduce a Test Me! button (or something similar) which lets Oppia automatically pick which skills the learner should practice.

When setting up a practice session, learners may specify one of two formats:
Answer a specific number of questions (each question must be answered correctly before continuing to the next one)
Reach a specified score (visible to the learner); this is the default selected option when starting a new practice session


'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 17: Score ---
$s = $p.Slides.Add(17, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Score'
$body = @'
A score is a value assigned to each question answered by the learner in a test. Scores are used to determine whether a learner has sufficiently demonstrated an understanding of the skills being tested.

Scores may be visible to the learner depending on the type of question set thats being shown. Scores will be shown as a progress bar with an optional animation when a score is achieved. Larger score numbers are preferred (e.g. 100 points for each correct answer without any help).

Scores may have different values depending on how much help the learner requires when answering the question. For example, a correct answer without any help may provide a score of 100. If the learner uses a hint, they may only get a score of 90 (where each hint reduces the score by 10). If the learner uses the solution in order to pass the question, they receive no score. Scores can never be negative.

'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 18: Question ---
$s = $p.Slides.Add(18, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Question'
$body = @'
A question is a prompt and interaction with a correct answer to evaluate the learners understanding of a specific skill. Question data comes from a configurable bank of questions and answers that allows Oppia to automatically generate a question in real-time. Oppia will be able to generate some data formulaically. Oppia will also provide hints and an optional solution to the question for cases when the learner is struggling.

If the learner is clearly struggling on a question (ie, they needed to ask for a solution to questions of the same topic several times), then they will be shown a concept card for the skill they are struggling on.
'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

# --- Slide 19: Skill ---
$s = $p.Slides.Add(19, 2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = 'Skill'
$body = @'
Skills are atomic units of learning that should represent a single concrete idea. Skills may occasionally need to be broken down further (e.g. to teach equivalent fractions, the concept of equivalence may need to be separately taught) using other skills. This means skills can build on one another: equivalence in the context of fractions could block learning about equivalent fractions, but these can be modeled as separate skills.

Note that cases where a skill seems to comprise many smaller parts may actually indicate a concept that should be modeled as a topic instead. Skills occasionally need prerequisite skills to help explain them, but conceptually should remain small, specific building blocks.


'@
$s.Shapes.Item(2).TextFrame.TextRange.Text = $body

